# The workbook lists image filenames used by the gallery. Two of the
# filenames in the sheet no longer match the actual asset names on disk,
# so correct them in place (same cells, same styling):
#   - "Still Life in White #2.jpg"  -> "Still-Life in White #2.jpg"  (row 9)
#   - "Redwood Picnic.jpeg"         -> "Redwood Picnic.jpg"          (row 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Still-Life in White #2.jpg"
$ws.Range("A13").Value = "Redwood Picnic.jpg"
